# Populate the assignee names in the "2025.1" sheet (column D) with the
# full names, added as new shared strings:
#   "Adel"  -> "Adel Jambalos"
#   "Ave"   -> "Ave Manriquez"
#   "Verna" -> "Verna Deatras"
# New strings must be appended in this (alphabetical) order so they land
# at shared-string indices 281, 282, 283 respectively, matching the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$adelCells = @("D15", "D19", "D36", "D37", "D51", "D64", "D69", "D71", "D72", "D87", "D117", "D118", "D121", "D122", "D135", "D137", "D139", "D140", "D143", "D169", "D175", "D179", "D182", "D187", "D209", "D213", "D217", "D224", "D228", "D230")

$aveCells = @("D16", "D24", "D73", "D120", "D134", "D136", "D142", "D146", "D153", "D155", "D158", "D159", "D160", "D168", "D172", "D174", "D178", "D185", "D188", "D208", "D212", "D216", "D221", "D227")

$vernaCells = @("D8", "D9", "D18", "D20", "D22", "D25", "D27", "D28", "D29", "D30", "D32", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D52", "D54", "D55", "D57", "D60", "D62", "D63", "D67", "D70", "D74", "D102", "D103", "D104", "D108", "D109", "D112", "D113", "D114", "D115", "D116", "D119", "D124", "D129", "D131", "D132", "D133", "D138", "D144", "D145", "D148", "D149", "D150", "D151", "D156", "D161", "D167", "D171", "D173", "D177", "D181", "D184", "D189", "D191", "D192", "D193", "D194", "D201", "D211", "D215", "D219", "D223", "D226", "D231")

foreach ($c in $adelCells) {
    $ws.Range($c).Value = "Adel Jambalos"
}

foreach ($c in $aveCells) {
    $ws.Range($c).Value = "Ave Manriquez"
}

foreach ($c in $vernaCells) {
    $ws.Range($c).Value = "Verna Deatras"
}

# Turn on the AutoFilter for the data range.
$ws.Range("A1:L231").AutoFilter()

# Restore the last active cell/selection recorded in the sheet view.
$ws.Range("D232").Select()
